$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 975.8570999999999
$ws.Range("I11").Value = 975.8570999999999
$ws.Range("K11").Value = 975.8570999999999
$ws.Range("M11").Value = -835.8570999999999
$ws.Range("H33").Value = 974
$ws.Range("I33").Value = 397.33334
$ws.Range("K33").Value = 397.33334
$ws.Range("M33").Value = -168.33334
$ws.Range("H43").Value = 5532.8887
$ws.Range("I43").Value = 5399.8
$ws.Range("J43").Value = 5699.25
$ws.Range("K43").Value = 5399.8
$ws.Range("L43").Value = 5699.25
$ws.Range("M43").Value = -5330.8
$ws.Range("N43").Value = -5837.25
$ws.Range("H64").Value = 22226278
$ws.Range("I64").Value = 25003938
$ws.Range("K64").Value = 25003938
$ws.Range("M64").Value = -25003690
$ws.Range("H67").Value = 22226278
$ws.Range("I67").Value = 25003938
$ws.Range("K67").Value = 25003938
$ws.Range("M67").Value = -25003080
$ws.Range("H69").Value = 19847.5
$ws.Range("I69").Value = 9390
$ws.Range("K69").Value = 28170
$ws.Range("M69").Value = -27296
$ws.Range("H72").Value = 19847.5
$ws.Range("I72").Value = 9390
$ws.Range("K72").Value = 84510
$ws.Range("M72").Value = -80142
$ws.Range("H74").Value = 55558040
$ws.Range("I74").Value = 55558040
$ws.Range("K74").Value = 55558040
$ws.Range("M74").Value = -55557104
$ws.Range("H77").Value = 55558040
$ws.Range("I77").Value = 55558040
$ws.Range("K77").Value = 277790200
$ws.Range("M77").Value = -277785520
$ws.Range("H138").Value = 1146.9584
$ws.Range("I138").Value = 887.9545000000001
$ws.Range("J138").Value = 3996
$ws.Range("K138").Value = 2663.8635
$ws.Range("L138").Value = 11988
$ws.Range("M138").Value = 2476.1365
$ws.Range("N138").Value = -22268

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 895.7143
$ws.Range("I4").Value = 895.7143
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 895.7143
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -779.7143
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 820.13043
$ws.Range("I5").Value = 715.3333
$ws.Range("K5").Value = 715.3333
$ws.Range("M5").Value = -603.3333
$ws.Range("H6").Value = 21577.8
$ws.Range("I6").Value = 21002
$ws.Range("J6").Value = 21721.75
$ws.Range("K6").Value = 21002
$ws.Range("L6").Value = 21721.75
$ws.Range("M6").Value = -20829
$ws.Range("N6").Value = -22067.75
$ws.Range("H74").Value = 32261826
$ws.Range("J74").Value = 5384
$ws.Range("L74").Value = 5384
$ws.Range("N74").Value = -7132
$ws.Range("H77").Value = 32261826
$ws.Range("J77").Value = 5384
$ws.Range("L77").Value = 26920
$ws.Range("N77").Value = -35656
$ws.Range("H132").Value = 3228479.5
$ws.Range("J132").Value = 800
$ws.Range("L132").Value = 2400
$ws.Range("N132").Value = -7460
$ws.Range("H139").Value = 149999.75
$ws.Range("J139").Value = 149999.75
$ws.Range("L139").Value = 149999.75
$ws.Range("N139").Value = -160279.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 820.13043
$ws.Range("I4").Value = 715.3333
$ws.Range("K4").Value = 715.3333
$ws.Range("M4").Value = -600.3333
$ws.Range("H22").Value = 1444765.1
$ws.Range("I22").Value = 1241.8235
$ws.Range("J22").Value = 2978508.5
$ws.Range("K22").Value = 1241.8235
$ws.Range("L22").Value = 2978508.5
$ws.Range("M22").Value = -1068.8235
$ws.Range("N22").Value = -2978854.5
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H134").Value = 42919660
$ws.Range("I134").Value = 51502492
$ws.Range("K134").Value = 154507476
$ws.Range("M134").Value = -154504941

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1135.1428
$ws.Range("I7").Value = 824.3333
$ws.Range("K7").Value = 824.3333
$ws.Range("M7").Value = -711.3333
$ws.Range("H22").Value = 9365.727999999999
$ws.Range("I22").Value = 20154.6
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 20154.6
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = -19804.6
$ws.Range("N22").Value = -1075
$ws.Range("H31").Value = 5305
$ws.Range("I31").Value = 3761.7144
$ws.Range("K31").Value = 3761.7144
$ws.Range("M31").Value = -3466.7144
$ws.Range("H34").Value = 5305
$ws.Range("I34").Value = 3761.7144
$ws.Range("K34").Value = 3761.7144
$ws.Range("M34").Value = -3559.7144
$ws.Range("H58").Value = 33342154
$ws.Range("I58").Value = 50011650
$ws.Range("K58").Value = 50011650
$ws.Range("M58").Value = -50011447
$ws.Range("H136").Value = 33342154
$ws.Range("I136").Value = 50011650
$ws.Range("K136").Value = 150034950
$ws.Range("M136").Value = -150032400

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 48588.855
$ws.Range("I5").Value = 77645.38
$ws.Range("J5").Value = 1372
$ws.Range("K5").Value = 232936.14
$ws.Range("L5").Value = 4116
$ws.Range("M5").Value = -232824.14
$ws.Range("N5").Value = -4340
$ws.Range("H6").Value = 178.5
$ws.Range("I6").Value = 178.5
$ws.Range("K6").Value = 535.5
$ws.Range("M6").Value = -422.5
$ws.Range("H37").Value = 132892.08
$ws.Range("J37").Value = 132892.08
$ws.Range("L37").Value = 398676.24
$ws.Range("N37").Value = -398900.24
$ws.Range("H135").Value = 48588.855
$ws.Range("I135").Value = 77645.38
$ws.Range("J135").Value = 1372
$ws.Range("K135").Value = 698808.42
$ws.Range("L135").Value = 12348
$ws.Range("M135").Value = -696273.42
$ws.Range("N135").Value = -17418

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 143.3158
$ws.Range("I2").Value = 82.30768999999999
$ws.Range("J2").Value = 275.5
$ws.Range("K2").Value = 82.30768999999999
$ws.Range("L2").Value = 275.5
$ws.Range("M2").Value = 30.69231000000001
$ws.Range("N2").Value = -501.5
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("M21").Value = -49827
$ws.Range("H30").Value = 50000
$ws.Range("I30").Value = 50000
$ws.Range("K30").Value = 50000
$ws.Range("M30").Value = -49895
$ws.Range("H35").Value = 24750
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 24750
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 24750
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -25346
$ws.Range("H36").Value = 6553.4
$ws.Range("I36").Value = 4908.5
$ws.Range("J36").Value = 7650
$ws.Range("K36").Value = 4908.5
$ws.Range("L36").Value = 7650
$ws.Range("M36").Value = -4423.5
$ws.Range("N36").Value = -8620
$ws.Range("H132").Value = 4169348
$ws.Range("I132").Value = 4810069
$ws.Range("K132").Value = 14430207
$ws.Range("M132").Value = -14427677

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 792.5278
$ws.Range("I93").Value = 802.14813
$ws.Range("K93").Value = 802.14813
$ws.Range("M93").Value = 445.85187
$ws.Range("H136").Value = 2484
$ws.Range("I136").Value = 1718.875
$ws.Range("J136").Value = 3249.125
$ws.Range("K136").Value = 5156.625
$ws.Range("L136").Value = 9747.375
$ws.Range("M136").Value = -2606.625
$ws.Range("N136").Value = -14847.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1045.3334
$ws.Range("I113").Value = 994.61536
$ws.Range("K113").Value = 2983.84608
$ws.Range("M113").Value = -813.8460800000003
$ws.Range("H136").Value = 7247956.5
$ws.Range("I136").Value = 7354441.5
$ws.Range("K136").Value = 22063324.5
$ws.Range("M136").Value = -22060774.5
